$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Row 2 stays same text, but underlying shared-string order changes due to
# new strings inserted mid-list (handled automatically by the engine) - just
# rewrite the same values to be safe.
$ws.Range("A2").Value = 19330051920441
$ws.Range("B2").Value = "GARCIA"
$ws.Range("C2").Value = "ANTONIO"
$ws.Range("D2").Value = "ABRAHAM"
$ws.Range("E2").Value = "INGLÉS IV"
$ws.Range("F2").Value = "4APV"
$ws.Range("G2").Value = 2

# New row 3
$ws.Range("A3").Value = 19330051920414
$ws.Range("B3").Value = "TREJO"
$ws.Range("C3").Value = "LUENGAS"
$ws.Range("D3").Value = "ELIZABETH"
$ws.Range("E3").Value = "INGLÉS IV"
$ws.Range("F3").Value = "4ASV"
$ws.Range("G3").Value = 2

# New row 4
$ws.Range("A4").Value = 19330051920053
$ws.Range("B4").Value = "ESPINOSA"
$ws.Range("C4").Value = "TZOPITL"
$ws.Range("D4").Value = "YASIEL"
$ws.Range("E4").Value = "INGLÉS IV"
$ws.Range("F4").Value = "4AEV"
$ws.Range("G4").Value = 1
